$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New school rows (18: McKee Road, 19: Providence Spring, 20: Polo Ridge) ---
# Shared strings must be introduced in this order so the new <si> entries in
# sharedStrings.xml land as: 30=Providence Spring, 31=McKee Road, 32=Polo Ridge.
$ws.Range("A19").Value = "Providence Spring"
$ws.Range("A18").Value = "McKee Road"
$ws.Range("A20").Value = "Polo Ridge"

# Match the existing centered-number style (style index 2) used by B2:N17.
$ws.Range("B18:N20").HorizontalAlignment = -4108

# Row 18 - McKee Road
$ws.Range("C18").Value = 86.4
$ws.Range("D18").Value = 86.1
$ws.Range("E18").Value = 80.4
$ws.Range("F18").Value = 91.8
$ws.Range("G18").Value = 88.4
$ws.Range("H18").Value = 514
$ws.Range("I18").Value = 15.8
$ws.Range("J18").Value = 66.3
$ws.Range("K18").Value = 7.8
$ws.Range("L18").Value = 8.6
$ws.Range("M18").Value = 13.4
$ws.Range("N18").Value = 3.9

# Row 19 - Providence Spring
$ws.Range("C19").Value = 92.7
$ws.Range("D19").Value = 80
$ws.Range("E19").Value = 91.2
$ws.Range("F19").Value = 93.1
$ws.Range("G19").Value = 95
$ws.Range("H19").Value = 932
$ws.Range("I19").Value = 18.9
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 4.2
$ws.Range("L19").Value = 5
$ws.Range("M19").Value = 18.6
$ws.Range("N19").Value = 2.3

# Row 20 - Polo Ridge
$ws.Range("C20").Value = 89
$ws.Range("D20").Value = 82.8
$ws.Range("E20").Value = 86.5
$ws.Range("F20").Value = 90
$ws.Range("G20").Value = 93.2
$ws.Range("H20").Value = 1021
$ws.Range("I20").Value = 17.8
$ws.Range("J20").Value = 53.5
$ws.Range("K20").Value = 7.1
$ws.Range("L20").Value = 5.3
$ws.Range("M20").Value = 31.1
$ws.Range("N20").Value = 3.1

# Performance Score = (0.8 * Achievement) + (0.2 * Growth), same formula used
# by the existing B3:B17 shared-formula column; extends it down to row 20.
$ws.Range("B18:B20").FormulaR1C1 = "=(0.8*RC[1])+(0.2*RC[2])"

# Reflect the updated selection recorded in the workbook.
$ws.Range("N22").Select()
